# Update: Actualización desde MV -datos-
# Appends the new quarterly row (01-07-2021) to the "Tasas de depósito en
# bolsa (mercado secundario)" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 60

# The date-like label "01-07-2021" is ambiguous as a date (both "01" and
# "07" are valid month numbers), so Excel's automatic type detection would
# otherwise convert it into a date serial number. Force the cell to Text
# first so it is stored as a plain string (matching the rest of column A),
# then drop the cell back to the default "Normal" style so no stray
# number-format is left applied to the cell.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "01-07-2021"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 0.9399999999999999
$ws.Cells.Item($row, 3).Value = 1.46
$ws.Cells.Item($row, 4).Value = 1.95
$ws.Cells.Item($row, 5).Value = 2.07
$ws.Cells.Item($row, 6).Value = -1.17
